$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"25.35940266666667"
$ws.Range("H2").Value = [double]"76.078208"
$ws.Range("I2").Value = [double]"0.005186643687654987"
$ws.Range("J2").Value = [double]"0.005186643687654986"
$ws.Range("M2").Value = [double]"0.3360566666666667"
$ws.Range("N2").Value = [double]"1.00817"
$ws.Range("O2").Value = [double]"0.01570866217798777"
$ws.Range("P2").Value = [double]"0.01570866217798777"
$ws.Range("Q2").Value = [double]"8.522196328817778"
$ws.Range("R2").Value = [double]"76.69976695936"
$ws.Range("S2").Value = [double]"8.147523352696493e-05"
$ws.Range("T2").Value = [double]"8.147523352696491e-05"

$ws.Range("G3").Value = [double]"25.35940266666667"
$ws.Range("H3").Value = [double]"76.078208"
$ws.Range("I3").Value = [double]"0.005186643687654987"
$ws.Range("J3").Value = [double]"0.005186643687654986"
$ws.Range("O3").Value = [double]"0.109316751024163"
$ws.Range("P3").Value = [double]"0.1093167510241629"
$ws.Range("Q3").Value = [double]"59.30605698312534"
$ws.Range("R3").Value = [double]"533.754512848128"
$ws.Range("S3").Value = [double]"0.0005669870366544266"
$ws.Range("T3").Value = [double]"0.0005669870366544265"

$ws.Range("G4").Value = [double]"25.35940266666667"
$ws.Range("H4").Value = [double]"76.078208"
$ws.Range("I4").Value = [double]"0.005186643687654987"
$ws.Range("J4").Value = [double]"0.005186643687654986"
$ws.Range("M4").Value = [double]"18.491866"
$ws.Range("N4").Value = [double]"55.47559800000001"
$ws.Range("O4").Value = [double]"0.864385399390831"
$ws.Range("P4").Value = [double]"0.864385399390831"
$ws.Range("Q4").Value = [double]"468.9426759520427"
$ws.Range("R4").Value = [double]"4220.484083568384"
$ws.Range("S4").Value = [double]"0.004483259075451588"
$ws.Range("T4").Value = [double]"0.004483259075451588"

$ws.Range("G5").Value = [double]"25.35940266666667"
$ws.Range("H5").Value = [double]"76.078208"
$ws.Range("I5").Value = [double]"0.005186643687654987"
$ws.Range("J5").Value = [double]"0.005186643687654986"
$ws.Range("M5").Value = [double]"0.2265353333333333"
$ws.Range("N5").Value = [double]"0.6796059999999999"
$ws.Range("O5").Value = [double]"0.01058918740701822"
$ws.Range("P5").Value = [double]"0.01058918740701822"
$ws.Range("Q5").Value = [double]"5.744800736227555"
$ws.Range("R5").Value = [double]"51.703206626048"
$ws.Range("S5").Value = [double]"5.492234202200672e-05"
$ws.Range("T5").Value = [double]"5.492234202200671e-05"

$ws.Range("I6").Value = [double]"0.9837462940761621"
$ws.Range("J6").Value = [double]"0.983746294076162"
$ws.Range("M6").Value = [double]"0.3360566666666667"
$ws.Range("N6").Value = [double]"1.00817"
$ws.Range("O6").Value = [double]"0.01570866217798777"
$ws.Range("P6").Value = [double]"0.01570866217798777"
$ws.Range("Q6").Value = [double]"1616.397724759542"
$ws.Range("R6").Value = [double]"14547.57952283588"
$ws.Range("S6").Value = [double]"0.01545333820248984"
$ws.Range("T6").Value = [double]"0.01545333820248984"

$ws.Range("I7").Value = [double]"0.9837462940761621"
$ws.Range("J7").Value = [double]"0.983746294076162"
$ws.Range("O7").Value = [double]"0.109316751024163"
$ws.Range("P7").Value = [double]"0.1093167510241629"
$ws.Range("S7").Value = [double]"0.1075399487004668"
$ws.Range("T7").Value = [double]"0.1075399487004668"

$ws.Range("I8").Value = [double]"0.9837462940761621"
$ws.Range("J8").Value = [double]"0.983746294076162"
$ws.Range("M8").Value = [double]"18.491866"
$ws.Range("N8").Value = [double]"55.47559800000001"
$ws.Range("O8").Value = [double]"0.864385399390831"
$ws.Range("P8").Value = [double]"0.864385399390831"
$ws.Range("Q8").Value = [double]"88943.95824798895"
$ws.Range("R8").Value = [double]"800495.6242319006"
$ws.Range("S8").Value = [double]"0.8503359333042733"
$ws.Range("T8").Value = [double]"0.8503359333042732"

$ws.Range("I9").Value = [double]"0.9837462940761621"
$ws.Range("J9").Value = [double]"0.983746294076162"
$ws.Range("M9").Value = [double]"0.2265353333333333"
$ws.Range("N9").Value = [double]"0.6796059999999999"
$ws.Range("O9").Value = [double]"0.01058918740701822"
$ws.Range("P9").Value = [double]"0.01058918740701822"
$ws.Range("Q9").Value = [double]"1089.61146645202"
$ws.Range("R9").Value = [double]"9806.503198068183"
$ws.Range("S9").Value = [double]"0.01041707386893214"
$ws.Range("T9").Value = [double]"0.01041707386893214"

$ws.Range("G10").Value = [double]"51.27300266666666"
$ws.Range("H10").Value = [double]"153.819008"
$ws.Range("I10").Value = [double]"0.01048663484403512"
$ws.Range("J10").Value = [double]"0.01048663484403512"
$ws.Range("M10").Value = [double]"0.3360566666666667"
$ws.Range("N10").Value = [double]"1.00817"
$ws.Range("O10").Value = [double]"0.01570866217798777"
$ws.Range("P10").Value = [double]"0.01570866217798777"
$ws.Range("Q10").Value = [double]"17.23063436615111"
$ws.Range("R10").Value = [double]"155.07570929536"
$ws.Range("S10").Value = [double]"0.0001647310041488633"
$ws.Range("T10").Value = [double]"0.0001647310041488633"

$ws.Range("G11").Value = [double]"51.27300266666666"
$ws.Range("H11").Value = [double]"153.819008"
$ws.Range("I11").Value = [double]"0.01048663484403512"
$ws.Range("J11").Value = [double]"0.01048663484403512"
$ws.Range("O11").Value = [double]"0.109316751024163"
$ws.Range("P11").Value = [double]"0.1093167510241629"
$ws.Range("Q11").Value = [double]"119.9081720423253"
$ws.Range("R11").Value = [double]"1079.173548380928"
$ws.Range("S11").Value = [double]"0.0011463648503267"
$ws.Range("T11").Value = [double]"0.001146364850326699"

$ws.Range("G12").Value = [double]"51.27300266666666"
$ws.Range("H12").Value = [double]"153.819008"
$ws.Range("I12").Value = [double]"0.01048663484403512"
$ws.Range("J12").Value = [double]"0.01048663484403512"
$ws.Range("M12").Value = [double]"18.491866"
$ws.Range("N12").Value = [double]"55.47559800000001"
$ws.Range("O12").Value = [double]"0.864385399390831"
$ws.Range("P12").Value = [double]"0.864385399390831"
$ws.Range("Q12").Value = [double]"948.1334947296427"
$ws.Range("R12").Value = [double]"8533.201452566784"
$ws.Range("S12").Value = [double]"0.009064494047927106"
$ws.Range("T12").Value = [double]"0.009064494047927106"

$ws.Range("G13").Value = [double]"51.27300266666666"
$ws.Range("H13").Value = [double]"153.819008"
$ws.Range("I13").Value = [double]"0.01048663484403512"
$ws.Range("J13").Value = [double]"0.01048663484403512"
$ws.Range("M13").Value = [double]"0.2265353333333333"
$ws.Range("N13").Value = [double]"0.6796059999999999"
$ws.Range("O13").Value = [double]"0.01058918740701822"
$ws.Range("P13").Value = [double]"0.01058918740701822"
$ws.Range("Q13").Value = [double]"11.61514675009422"
$ws.Range("R13").Value = [double]"104.536320750848"
$ws.Range("S13").Value = [double]"0.0001110449416324552"
$ws.Range("T13").Value = [double]"0.0001110449416324552"

$ws.Range("G14").Value = [double]"2.837922333333333"
$ws.Range("H14").Value = [double]"8.513767"
$ws.Range("I14").Value = [double]"0.0005804273921477663"
$ws.Range("J14").Value = [double]"0.0005804273921477662"
$ws.Range("M14").Value = [double]"0.3360566666666667"
$ws.Range("N14").Value = [double]"1.00817"
$ws.Range("O14").Value = [double]"0.01570866217798777"
$ws.Range("P14").Value = [double]"0.01570866217798777"
$ws.Range("Q14").Value = [double]"0.9537027195988889"
$ws.Range("R14").Value = [double]"8.58332447639"
$ws.Range("S14").Value = [double]"9.117737822099693e-06"
$ws.Range("T14").Value = [double]"9.117737822099691e-06"

$ws.Range("G15").Value = [double]"2.837922333333333"
$ws.Range("H15").Value = [double]"8.513767"
$ws.Range("I15").Value = [double]"0.0005804273921477663"
$ws.Range("J15").Value = [double]"0.0005804273921477662"
$ws.Range("O15").Value = [double]"0.109316751024163"
$ws.Range("P15").Value = [double]"0.1093167510241629"
$ws.Range("Q15").Value = [double]"6.636827603024666"
$ws.Range("R15").Value = [double]"59.731448427222"
$ws.Range("S15").Value = [double]"6.345043671502156e-05"
$ws.Range("T15").Value = [double]"6.345043671502154e-05"

$ws.Range("G16").Value = [double]"2.837922333333333"
$ws.Range("H16").Value = [double]"8.513767"
$ws.Range("I16").Value = [double]"0.0005804273921477663"
$ws.Range("J16").Value = [double]"0.0005804273921477662"
$ws.Range("M16").Value = [double]"18.491866"
$ws.Range("N16").Value = [double]"55.47559800000001"
$ws.Range("O16").Value = [double]"0.864385399390831"
$ws.Range("P16").Value = [double]"0.864385399390831"
$ws.Range("Q16").Value = [double]"52.47847950640734"
$ws.Range("R16").Value = [double]"472.306315557666"
$ws.Range("S16").Value = [double]"0.0005017129631790254"
$ws.Range("T16").Value = [double]"0.0005017129631790254"

$ws.Range("G17").Value = [double]"2.837922333333333"
$ws.Range("H17").Value = [double]"8.513767"
$ws.Range("I17").Value = [double]"0.0005804273921477663"
$ws.Range("J17").Value = [double]"0.0005804273921477662"
$ws.Range("M17").Value = [double]"0.2265353333333333"
$ws.Range("N17").Value = [double]"0.6796059999999999"
$ws.Range("O17").Value = [double]"0.01058918740701822"
$ws.Range("P17").Value = [double]"0.01058918740701822"
$ws.Range("Q17").Value = [double]"0.6428896817557777"
$ws.Range("R17").Value = [double]"5.786007135802"
$ws.Range("S17").Value = [double]"6.146254431619552e-06"
$ws.Range("T17").Value = [double]"6.146254431619551e-06"

